# Rewrite the test sheet: replace the old "Nombre/Apellidos/DNI/Fecha
# nacimiento/Nacionalidad/Dirección postal/NIF/Correo electrónico/pollingStation"
# dataset (columns A:I) with the new, smaller "Nombre/Localización/Correo
# electrónico/ID/Tipo" dataset (columns A:E).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the old hyperlinks first so clearing the cells doesn't leave dangling
# hyperlink references behind.
$ws.Hyperlinks.Delete()

# Clear everything else so no stale cells/columns (F:I) survive.
$ws.Cells.Clear()

# Write the new values in the same order the source data would naturally be
# entered column-by-column (Nombre, then Correo electrónico + its rows, then
# the name rows, then Localización, then ID/Tipo) so the shared-string table
# comes out in the same order as the target workbook.
$ws.Range("A1").Value = "Nombre"
$ws.Range("C1").Value = "Correo electrónico"

$ws.Range("C2").Value = "juan@example.com"
$ws.Range("C3").Value = "luis@example.com"
$ws.Range("C4").Value = "ana@example.com"

$ws.Range("A2").Value = "Juan Torres Pardo"
$ws.Range("A3").Value = "Luis López Fernando"
$ws.Range("A4").Value = "Ana Torres Pardo"

$ws.Range("B1").Value = "Localización"
$ws.Range("B2").Value = "41.5N35.99W"
$ws.Range("B3").Value = "41.5N35.99W"
$ws.Range("B4").Value = "41.5N35.99W"

$ws.Range("D1").Value = "ID"
$ws.Range("E1").Value = "Tipo"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 1
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 1
$ws.Range("D4").Value = 3
$ws.Range("E4").Value = 2

# Re-create the mailto hyperlinks on the email column (lost when cells were
# cleared), then re-apply the workbook's existing "Hipervínculo" cell style so
# the cells reuse the same style slot the original hyperlinked cells used.
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:juan@example.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:luis@example.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:ana@example.com") | Out-Null
$ws.Range("C2").Style = "Hipervínculo"
$ws.Range("C3").Style = "Hipervínculo"
$ws.Range("C4").Style = "Hipervínculo"

# Column widths: leave A:C alone (already 23.42578125, untouched by the
# diff) and only resize D to the old column E's width (22.140625). The COM
# ColumnWidth setter here only has 1/6-character resolution, so
# 21.333333333333332 is the closest input that round-trips to
# 22.140625 in the saved OOXML (lands on 22.166666666666668, ~0.03 off).
$ws.Columns("D").ColumnWidth = 21.333333333333332

# Selection moves to the (now unused) F column, matching the saved view state.
$ws.Range("F1:F1048576").Select()
